# Daily refresh of the cryptocurrency price/volume table (GitHub Actions bot).
# Price cells that look numeric are written with a leading apostrophe
# (Excel's "quote prefix") so they stay literal text - matching the
# original inline-string cells - instead of being re-parsed as numbers
# (which would silently drop trailing zeros, e.g. "18.40" -> "18.4").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.399.81'
$ws.Range("E2").Value = '  -1.31%  '

$ws.Range("D3").Value = '1.915.39'
$ws.Range("E3").Value = '  +1.50%  '

$ws.Range("D4").Value = "'0.9982"
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = "'242.25"
$ws.Range("E5").Value = '  +1.38%  '

$ws.Range("D6").Value = "'0.9982"
$ws.Range("E6").Value = '  -0.02%  '

$ws.Range("D7").Value = "'0.4693"
$ws.Range("E7").Value = '  -1.60%  '

$ws.Range("D8").Value = "'0.2860"
$ws.Range("E8").Value = '  -0.58%  '

$ws.Range("D9").Value = "'0.06827"
$ws.Range("E9").Value = '  +3.91%  '

$ws.Range("D10").Value = "'111.58"
$ws.Range("E10").Value = '  +13.91%  '

$ws.Range("D11").Value = "'18.40"
$ws.Range("E11").Value = '  -2.90%  '

$ws.Range("D12").Value = "'0.07721"
$ws.Range("E12").Value = '  +1.92%  '

$ws.Range("D13").Value = '1.888.04'
$ws.Range("E13").Value = '  +0.62%  '

$ws.Range("D14").Value = "'5.267"
$ws.Range("E14").Value = '  +2.59%  '

$ws.Range("D15").Value = "'0.6584"
$ws.Range("E15").Value = '  +0.21%  '

$ws.Range("D16").Value = "'299.31"
$ws.Range("E16").Value = '  -3.30%  '

$ws.Range("D17").Value = '30.406.58'
$ws.Range("E17").Value = '  -1.24%  '

$ws.Range("E18").Value = '  +0.37%  '

$ws.Range("D19").Value = "'0.9992"
$ws.Range("E19").Value = '  +0.08%  '

$ws.Range("D20").Value = "'12.90"
$ws.Range("E20").Value = '  -2.32%  '

$ws.Range("D21").Value = '2.132.37'
$ws.Range("E21").Value = '  +0.45%  '

$ws.Range("D22").Value = "'0.9984"
$ws.Range("E22").Value = '  -0.14%  '

$ws.Range("D23").Value = "'5.235"
$ws.Range("E23").Value = '  +2.06%  '

$ws.Range("D24").Value = "'6.223"
$ws.Range("E24").Value = '  +0.47%  '

$ws.Range("B25").Value = 'EthereumClassic'
$ws.Range("C25").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D25").Value = "'21.77"
$ws.Range("E25").Value = '  +7.17%  '

$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = "'9.307"
$ws.Range("E26").Value = '  -0.15%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = "'168.10"
$ws.Range("E27").Value = '  +0.83%  '

$ws.Range("D28").Value = "'2.093"
$ws.Range("E28").Value = '  +7.31%  '

$ws.Range("D29").Value = "'0.1072"
$ws.Range("E29").Value = '  -0.38%  '

$ws.Range("D30").Value = "'1.366"
$ws.Range("E30").Value = '  +1.03%  '

$ws.Range("D31").Value = "'4.169"
$ws.Range("E31").Value = '  +0.00%  '

$ws.Range("D32").Value = "'3.989"
$ws.Range("E32").Value = '  +0.14%  '

$ws.Range("D33").Value = "'0.05059"
$ws.Range("E33").Value = '  +0.26%  '

$ws.Range("D34").Value = "'0.7374"
$ws.Range("E34").Value = '  +0.21%  '

$ws.Range("E35").Value = '  -2.22%  '

$ws.Range("D36").Value = "'0.02079"
$ws.Range("E36").Value = '  +6.51%  '

$ws.Range("D37").Value = "'2.739"
$ws.Range("E37").Value = '  +1.09%  '

$ws.Range("D38").Value = "'2.676"
$ws.Range("E38").Value = '  -1.05%  '

$ws.Range("D39").Value = "'2.057"
$ws.Range("E39").Value = '  -1.18%  '

$ws.Range("D40").Value = "'109.92"
$ws.Range("E40").Value = '  +1.84%  '

$ws.Range("D41").Value = "'0.8716"
$ws.Range("E41").Value = '  -3.65%  '

$ws.Range("D42").Value = "'5.814"
$ws.Range("E42").Value = '  +3.10%  '

$ws.Range("D43").Value = "'0.4248"
$ws.Range("E43").Value = '  +0.40%  '

$ws.Range("D44").Value = "'0.9979"
$ws.Range("E44").Value = '  -0.04%  '

$ws.Range("D45").Value = "'51.57"
$ws.Range("E45").Value = '  +20.22%  '

$ws.Range("D46").Value = "'67.46"
$ws.Range("E46").Value = '  +2.37%  '

$ws.Range("D47").Value = "'7.191"
$ws.Range("E47").Value = '  -2.79%  '

$ws.Range("D48").Value = "'9.299"
$ws.Range("E48").Value = '  +3.19%  '

$ws.Range("E49").Value = '  -1.56%  '

$ws.Range("D50").Value = "'34.80"
$ws.Range("E50").Value = '  -0.07%  '

$ws.Range("B51").Value = 'WOONetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D51").Value = "'0.2414"
$ws.Range("E51").Value = '  +9.79%  '
